$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Font color fix: the numeric data columns (B:D) font color goes from
# theme color 1 to an explicit black RGB color.
$ws.Range("B2:D13").Font.Color = 0

# Row height adjustments for data rows 2-13 (18.75 -> 19.5 points)
$ws.Rows("2:13").RowHeight = 19.5

# Data adjustments (Capacite / column B)
$ws.Range("B2").Value = 9000
$ws.Range("B3").Value = 4000
$ws.Range("B4").Value = 1400
$ws.Range("B5").Value = 1400
$ws.Range("B6").Value = 1400
$ws.Range("B7").Value = 1400
$ws.Range("B8").Value = 1400
$ws.Range("B9").Value = 1400
$ws.Range("B10").Value = 1400
$ws.Range("B11").Value = 1600
$ws.Range("B12").Value = 1600
$ws.Range("B13").Value = 1600
